$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAKE-OFF")

$ws.Range("C3").Value = 282.1452529462797
$ws.Range("C4").Value = 151.46515118137688
$ws.Range("C5").Value = 1427.3267887915854
$ws.Range("C6").Value = 1641.4258071103227
$ws.Range("C7").Value = 1853.2236402716744

$ws.Range("C10").Value = 925.6734020547234
$ws.Range("C11").Value = 496.93291070005546
$ws.Range("C12").Value = 4682.830671888403
$ws.Range("C13").Value = 5385.255272671663
$ws.Range("C14").Value = 6080.130053384759

$ws.Range("C17").Value = 63.5085873574696
$ws.Range("C18").Value = 68.3549361418386
$ws.Range("C20").Value = 74.326259768459
$ws.Range("C21").Value = 77.53034686818486

$ws.Range("C24").Value = 123.45081775750032
$ws.Range("C25").Value = 132.87136615044224
$ws.Range("C27").Value = 144.47869069462874
$ws.Range("C28").Value = 150.70693775673084

$ws.Range("C30").Value = 0.9966966966966815
$ws.Range("C31").Value = 1.0727547547547536
$ws.Range("C33").Value = 1.1664680426926468
$ws.Range("C34").Value = 1.2167526287794068

$ws.Range("C38").Value = 50.73686959720092

$ws.Range("C40").Value = 1180.8149251931632
$ws.Range("C44").Value = 172623.83873822342
$ws.Range("C46").Value = 63421.08699650115
